$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rebuild the lease-terms header row / data columns using pure cell-level
# moves (copy value+format from source range to destination range) instead
# of structural column insert/delete, to keep the worksheet's column
# metadata clean.
#
# Starting layout (A1:I3):
#   A Lease Reference | B Lease External Reference |
#   C Start Date Previous (empty) | D Value Previous (empty) |
#   E Start Date Current | F Value Current | G Start Date | H Value | I End Date
#
# Target layout (A1:H3):
#   A Lease Reference | B Lease External Reference |
#   C Start Date Previous Year | D End Date Previous Year | E Value Previous Year |
#   F Start Date | G End Date | H Value
#
# Column C picks up the old "Start Date Current" data (E), D is brand new
# "End Date Previous Year" data, E picks up the old "Value Current" data
# (F), F keeps "Start Date" (from G), G becomes "End Date" (from I), and H
# keeps "Value" in place. Each source column is read before it is
# overwritten, in left-to-right order, so no scratch cells are needed
# except column I which is read twice (for D's format, and for G) before
# finally being cleared.
# ---------------------------------------------------------------------------

$ws.Range("E1:E3").Copy($ws.Range("C1")) | Out-Null
$ws.Range("I1:I3").Copy($ws.Range("D1")) | Out-Null
$ws.Range("F1:F3").Copy($ws.Range("E1")) | Out-Null
$ws.Range("G1:G3").Copy($ws.Range("F1")) | Out-Null
$ws.Range("I1:I3").Copy($ws.Range("G1")) | Out-Null
$ws.Range("I1:I3").Clear() | Out-Null

# ---------------------------------------------------------------------------
# New "End Date Previous Year" data (column D): 2010-12-31 for both leases,
# keeping the date format that was copied in from column I above.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value2 = 40543
$ws.Range("D3").Value2 = 40543

# ---------------------------------------------------------------------------
# Header row text.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value2 = "Start Date Previous Year"
$ws.Range("D1").Value2 = "End Date Previous Year"
$ws.Range("E1").Value2 = "Value Previous Year"
$ws.Range("F1").Value2 = "Start Date"
$ws.Range("G1").Value2 = "End Date"
$ws.Range("H1").Value2 = "Value"

# ---------------------------------------------------------------------------
# Column widths: A and B stay as-is; C:D are widened to fit the new long
# "... Previous Year" headers; E:H keep their (shifted) best-fit widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 22.166666666666664
$ws.Columns.Item(4).ColumnWidth = 22.166666666666664

# ---------------------------------------------------------------------------
# Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("E1").Select() | Out-Null
